# Auto commit at 2025-11-07 11:50:27.39
# Update the Metrics sheet source values; dependent formulas on other
# sheets (today!B11:B22, E11:E22, F11:F22, and the volatile TODAY()-1 in
# today!A1) will recalculate automatically.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 77444.689999999988
$metrics.Range("B3").Value = 67632.2
$metrics.Range("B4").Value = 24009.7
$metrics.Range("B5").Value = 3260
$metrics.Range("B6").Value = 4873690.4399999995
$metrics.Range("B7").Value = 4109708.8800000008
$metrics.Range("B8").Value = 1430969.5299999998
$metrics.Range("B9").Value = 189467
$metrics.Range("B10").Value = 33339071.430000003
$metrics.Range("B11").Value = 31384984.039999999
$metrics.Range("B12").Value = 11712691.570000002
$metrics.Range("B13").Value = 1287097

# Move the active selection on the Metrics sheet to match the saved view.
$metrics.Activate()
$metrics.Range("E16").Select()

# Recalculate the whole workbook so dependent sheets (today, etc.) pick up
# the new values as well as the volatile TODAY()-1 formula.
$excel.Calculate()

# "today" is the sheet that was active/selected in the original workbook;
# restore it as the active sheet (with its own updated selection) so the
# saved view matches the source workbook.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F6").Select()
